$wb = $excel.ActiveWorkbook

# --- Task A1 ---
$wsA1 = $wb.Worksheets.Item("A1")
$wsA1.Range("A2").Value = "F2026442046CAD4820BA81C7DD1F3355C074C5F42019D271204EFB93D8EF45F8"
$wsA1.Range("B2").Value = "vivens"

# --- Task A2 (two rows entered, column by column: A then B then C) ---
$wsA2 = $wb.Worksheets.Item("A2")
$wsA2.Range("A2").Value = "84A0A1D0574F35248587D2D899913C37E1071636461135F56ADFE95A747F57BA"
$wsA2.Range("A3").Value = "4ACF9F587EE44EBE1571E5CA8E96AE1F0AC03093A2DD4E798E4D4EF6A04F74BC"
$wsA2.Range("B2").Value = "vivens"
$wsA2.Range("B2").Copy($wsA2.Range("B3"))
$wsA2.Range("C2").Value = "nft0001"
$wsA2.Range("C3").Value = "nft0002"

# --- Task A3 ---
$wsA3 = $wb.Worksheets.Item("A3")
$wsA3.Range("A2").Value = "AA08F2AB49F70B17AF67E4E09C9647B1C4BF11BD06520C577F0164A0B5F6A1BF"
$wsA3.Range("B2").Value = "juno1at6nu0jt7lzv0537mavw6k65kn3ekzv8lxmevfgvw3r2dfdshhpq9hwl5s"
$wsA3.Range("C2").Value = "nft0001"
$wsA3.Range("D2").Value = "uni-6"

# --- Task A5 (filled before A4, matching original authoring order) ---
$wsA5 = $wb.Worksheets.Item("A5")
$wsA5.Range("A2").Value = "F0E335D85081F77F0E249CB2157A0F20ABBB851D3DFC2D373F8F3F9A843A7525"
$wsA5.Range("B2").Value = "juno1at6nu0jt7lzv0537mavw6k65kn3ekzv8lxmevfgvw3r2dfdshhpq9hwl5s"
$wsA5.Range("C2").Value = "nft0001"
$wsA5.Range("D2").Value = "uni-6"

# --- Task A4 ---
$wsA4 = $wb.Worksheets.Item("A4")
$wsA4.Range("A2").Value = "E8659583AD6C648D0BFBA91B2026895F2B18C3A2B2F3E5FCE5D31FCDB5D4D9B7"
$wsA4.Range("B2").Value = "ibc/39158EBE0DF416D1D6C3E7F99C2CD113F04E59C7BAC8B5CCEF7563C358D898AE"
$wsA4.Range("C2").Value = "nft0002"
$wsA4.Range("D2").Value = "gon-flixnet-1"

# --- Task A6 ---
$wsA6 = $wb.Worksheets.Item("A6")
$wsA6.Range("A2").Value = "D76BDBEE44288E78DB03BF0100E4118A7488F1967A835DF051CD10991697CD4E"
$wsA6.Range("B2").Value = "ibc/39158EBE0DF416D1D6C3E7F99C2CD113F04E59C7BAC8B5CCEF7563C358D898AE"
$wsA6.Range("C2").Value = "nft0002"
$wsA6.Range("D2").Value = "gon-flixnet-1"

# --- View / selection state ---
$wsInfo = $wb.Worksheets.Item("Info")
$wsInfo.Range("E8").Select()

$wsA1.Range("B2").Select()
$wsA2.Range("C4").Select()
$wsA3.Range("B2:D2").Select()
$wsA4.Range("B2:D2").Select()
$wsA5.Range("B11").Select()

$wsA6.Activate()
$wsA6.Range("C27").Select()
